$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates in D4:D5 (44574 -> 44559) and D6:D7 (44559 -> 44574)
$ws.Range("D4").Value = 44559
$ws.Range("D5").Value = 44559
$ws.Range("D6").Value = 44574
$ws.Range("D7").Value = 44574
